$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 65 ---
$ws.Cells.Item(65, 4).Value = 44474     # D65 (Fecha)
$ws.Cells.Item(65, 9).Value = "Segunda" # I65 (Calidad)
$ws.Cells.Item(65, 10).Value = 500      # J65 (Volumen)
$ws.Cells.Item(65, 11).Value = 600      # K65 (Precio minimo)
$ws.Cells.Item(65, 12).Value = 700      # L65 (Precio maximo)
$ws.Cells.Item(65, 13).Value = 650      # M65 (Precio promedio ponderado)
$ws.Cells.Item(65, 16).Value = 650      # P65 (Precio $/Kg)

# --- Update existing row 66 ---
$ws.Cells.Item(66, 4).Value = 44474     # D66
$ws.Cells.Item(66, 9).Value = "Tercera" # I66
$ws.Cells.Item(66, 10).Value = 1200     # J66
$ws.Cells.Item(66, 11).Value = 400      # K66
$ws.Cells.Item(66, 12).Value = 500      # L66
$ws.Cells.Item(66, 13).Value = 450      # M66
$ws.Cells.Item(66, 16).Value = 450      # P66

# --- Update existing row 67 ---
$ws.Cells.Item(67, 4).Value = 44411     # D67
$ws.Cells.Item(67, 9).Value = "Tercera" # I67
$ws.Cells.Item(67, 10).Value = 1500     # J67
$ws.Cells.Item(67, 11).Value = 500      # K67
$ws.Cells.Item(67, 12).Value = 600      # L67
$ws.Cells.Item(67, 13).Value = 550      # M67
$ws.Cells.Item(67, 16).Value = 550      # P67

# --- Update existing row 68 (D68 unchanged) ---
$ws.Cells.Item(68, 9).Value = "Primera" # I68
$ws.Cells.Item(68, 10).Value = 500      # J68
$ws.Cells.Item(68, 11).Value = 700      # K68
$ws.Cells.Item(68, 12).Value = 800      # L68
$ws.Cells.Item(68, 13).Value = 750      # M68
$ws.Cells.Item(68, 16).Value = 750      # P68

# --- Update existing row 69 (I69, L69 unchanged) ---
$ws.Cells.Item(69, 4).Value = 44292     # D69
$ws.Cells.Item(69, 10).Value = 1000     # J69
$ws.Cells.Item(69, 11).Value = 600      # K69
$ws.Cells.Item(69, 13).Value = 650      # M69
$ws.Cells.Item(69, 16).Value = 650      # P69

# --- Update existing row 70 (becomes a "Tercera" row dated 44292) ---
$ws.Cells.Item(70, 4).Value = 44292     # D70
$ws.Cells.Item(70, 9).Value = "Tercera" # I70
$ws.Cells.Item(70, 10).Value = 800      # J70
$ws.Cells.Item(70, 11).Value = 400      # K70
$ws.Cells.Item(70, 12).Value = 500      # L70
$ws.Cells.Item(70, 13).Value = 450      # M70
$ws.Cells.Item(70, 16).Value = 450      # P70

# --- Add new row 71 ---
$ws.Cells.Item(71, 1).Value = 1
$ws.Cells.Item(71, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(71, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(71, 4).Value = 44323
$ws.Cells.Item(71, 5).Value = 15
$ws.Cells.Item(71, 6).Value = 100112008
$ws.Cells.Item(71, 7).Value = "Coliflor"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Segunda"
$ws.Cells.Item(71, 10).Value = 800
$ws.Cells.Item(71, 11).Value = 650
$ws.Cells.Item(71, 12).Value = 700
$ws.Cells.Item(71, 13).Value = 675
$ws.Cells.Item(71, 14).Value = "$/unidad"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 675
$ws.Cells.Item(71, 17).Value = 1
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# --- Add new row 72 ---
$ws.Cells.Item(72, 1).Value = 1
$ws.Cells.Item(72, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value = 44323
$ws.Cells.Item(72, 5).Value = 15
$ws.Cells.Item(72, 6).Value = 100112008
$ws.Cells.Item(72, 7).Value = "Coliflor"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Tercera"
$ws.Cells.Item(72, 10).Value = 1500
$ws.Cells.Item(72, 11).Value = 500
$ws.Cells.Item(72, 12).Value = 600
$ws.Cells.Item(72, 13).Value = 550
$ws.Cells.Item(72, 14).Value = "$/unidad"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 550
$ws.Cells.Item(72, 17).Value = 1
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# --- Apply the date number format to the new D71/D72 cells (style used by column D elsewhere) ---
$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
$ws.Cells.Item(72, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
